$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 280091.78
$ws.Range("I17").Value2 = 1699
$ws.Range("J17").Value2 = 296467.8
$ws.Range("K17").Value2 = 5097
$ws.Range("L17").Value2 = 889403.3999999999
$ws.Range("M17").Value2 = -4929
$ws.Range("N17").Value2 = -889739.3999999999
$ws.Range("H40").Value2 = 2678.28
$ws.Range("J40").Value2 = 1365
$ws.Range("L40").Value2 = 1365
$ws.Range("N40").Value2 = -1715
$ws.Range("H62").Value2 = 75005050
$ws.Range("J62").Value2 = 83338420
$ws.Range("L62").Value2 = 83338420
$ws.Range("N62").Value2 = -83339668
$ws.Range("H65").Value2 = 75005050
$ws.Range("J65").Value2 = 83338420
$ws.Range("L65").Value2 = 416692100
$ws.Range("N65").Value2 = -416698340
$ws.Range("H113").Value2 = 6980.5
$ws.Range("I113").Value2 = 7890.1763
$ws.Range("K113").Value2 = 7890.1763
$ws.Range("M113").Value2 = -4636.1763
$ws.Range("H116").Value2 = 17167.77
$ws.Range("J116").Value2 = 5998.5
$ws.Range("L116").Value2 = 5998.5
$ws.Range("N116").Value2 = -12882.5
$ws.Range("H137").Value2 = 7780475.5
$ws.Range("I137").Value2 = 11236722
$ws.Range("J137").Value2 = 3920.25
$ws.Range("K137").Value2 = 33710166
$ws.Range("L137").Value2 = 11760.75
$ws.Range("M137").Value2 = -33707616
$ws.Range("N137").Value2 = -16860.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value2 = 851
$ws.Range("I4").Value2 = 800
$ws.Range("J4").Value2 = 902
$ws.Range("K4").Value2 = 800
$ws.Range("L4").Value2 = 902
$ws.Range("M4").Value2 = -684
$ws.Range("N4").Value2 = -1134
$ws.Range("H32").Value2 = 7289.6865
$ws.Range("I32").Value2 = 7289.6865
$ws.Range("J32").Value2 = 0
$ws.Range("K32").Value2 = 7289.6865
$ws.Range("L32").Value2 = 0
$ws.Range("N32").Value2 = -7002.6865
$ws.Range("H45").Value2 = 1119.5
$ws.Range("I45").Value2 = 1431.5
$ws.Range("J45").Value2 = 963.5
$ws.Range("K45").Value2 = 1431.5
$ws.Range("L45").Value2 = 963.5
$ws.Range("M45").Value2 = -1054.5
$ws.Range("N45").Value2 = -1717.5
$ws.Range("H74").Value2 = 111577.22
$ws.Range("I74").Value2 = 111577.22
$ws.Range("K74").Value2 = 111577.22
$ws.Range("M74").Value2 = -110703.22
$ws.Range("H77").Value2 = 111577.22
$ws.Range("I77").Value2 = 111577.22
$ws.Range("K77").Value2 = 557886.1
$ws.Range("M77").Value2 = -553518.1
$ws.Range("H110").Value2 = 10607.167
$ws.Range("I110").Value2 = 10162
$ws.Range("K110").Value2 = 10162
$ws.Range("M110").Value2 = -8117

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value2 = 285.5
$ws.Range("I22").Value2 = 285.5
$ws.Range("J22").Value2 = 0
$ws.Range("K22").Value2 = 285.5
$ws.Range("L22").Value2 = 0
$ws.Range("N22").Value2 = -112.5
$ws.Range("H105").Value2 = 2484
$ws.Range("I105").Value2 = 2484
$ws.Range("K105").Value2 = 2484
$ws.Range("M105").Value2 = -737
$ws.Range("H132").Value2 = 109998.4
$ws.Range("I132").Value2 = 90000
$ws.Range("J132").Value2 = 114998
$ws.Range("K132").Value2 = 90000
$ws.Range("L132").Value2 = 114998
$ws.Range("M132").Value2 = -84940
$ws.Range("N132").Value2 = -125118
$ws.Range("H134").Value2 = 1763.3269
$ws.Range("I134").Value2 = 1336.909
$ws.Range("K134").Value2 = 4010.727
$ws.Range("M134").Value2 = -1475.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 2180.8
$ws.Range("I16").Value2 = 2180.8
$ws.Range("J16").Value2 = 0
$ws.Range("K16").Value2 = 2180.8
$ws.Range("L16").Value2 = 0
$ws.Range("N16").Value2 = -1893.8
$ws.Range("H113").Value2 = 2180.8
$ws.Range("I113").Value2 = 2180.8
$ws.Range("J113").Value2 = 0
$ws.Range("K113").Value2 = 2180.8
$ws.Range("L113").Value2 = 0
$ws.Range("N113").Value2 = -10.80000000000018

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value2 = 1551.9574
$ws.Range("I39").Value2 = 899.1
$ws.Range("J39").Value2 = 1728.4054
$ws.Range("K39").Value2 = 2697.3
$ws.Range("L39").Value2 = 5185.216200000001
$ws.Range("M39").Value2 = -2403.3
$ws.Range("N39").Value2 = -5773.216200000001
$ws.Range("H55").Value2 = 5399.7856
$ws.Range("J55").Value2 = 7117
$ws.Range("L55").Value2 = 21351
$ws.Range("N55").Value2 = -21705
$ws.Range("H86").Value2 = 383.91666
$ws.Range("I86").Value2 = 405.27274
$ws.Range("J86").Value2 = 149
$ws.Range("K86").Value2 = 1215.81822
$ws.Range("L86").Value2 = 447
$ws.Range("M86").Value2 = -29.81822000000011
$ws.Range("N86").Value2 = -2819
$ws.Range("H89").Value2 = 383.91666
$ws.Range("I89").Value2 = 405.27274
$ws.Range("J89").Value2 = 149
$ws.Range("K89").Value2 = 3647.45466
$ws.Range("L89").Value2 = 1341
$ws.Range("M89").Value2 = 2280.54534
$ws.Range("N89").Value2 = -13197
$ws.Range("H107").Value2 = 492.3913
$ws.Range("J107").Value2 = 517.6667
$ws.Range("L107").Value2 = 1553.0001
$ws.Range("N107").Value2 = -5393.0001
$ws.Range("H121").Value2 = 33334224
$ws.Range("I121").Value2 = 41667788
$ws.Range("J121").Value2 = 23810152
$ws.Range("K121").Value2 = 125003364
$ws.Range("L121").Value2 = 71430456
$ws.Range("M121").Value2 = -125002054
$ws.Range("N121").Value2 = -71433076
$ws.Range("H122").Value2 = 1958.1666
$ws.Range("I122").Value2 = 1952
$ws.Range("J122").Value2 = 1961.25
$ws.Range("K122").Value2 = 17568
$ws.Range("L122").Value2 = 17651.25
$ws.Range("M122").Value2 = -15118
$ws.Range("N122").Value2 = -22551.25
$ws.Range("H132").Value2 = 5081.839
$ws.Range("I132").Value2 = 5914.6
$ws.Range("J132").Value2 = 1612
$ws.Range("K132").Value2 = 53231.4
$ws.Range("L132").Value2 = 14508
$ws.Range("M132").Value2 = -50701.4
$ws.Range("N132").Value2 = -19568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value2 = 12665.667
$ws.Range("I10").Value2 = 13999.5
$ws.Range("J10").Value2 = 9998
$ws.Range("K10").Value2 = 13999.5
$ws.Range("L10").Value2 = 9998
$ws.Range("M10").Value2 = -13830.5
$ws.Range("N10").Value2 = -10336
$ws.Range("H102").Value2 = 1978.4736
$ws.Range("I102").Value2 = 1978.4736
$ws.Range("K102").Value2 = 1978.4736
$ws.Range("M102").Value2 = -356.4736
$ws.Range("H122").Value2 = 2297.2727
$ws.Range("I122").Value2 = 2391.3333
$ws.Range("K122").Value2 = 7173.999899999999
$ws.Range("M122").Value2 = -4723.999899999999
$ws.Range("H132").Value2 = 34198.145
$ws.Range("I132").Value2 = 43357.707
$ws.Range("K132").Value2 = 130073.121
$ws.Range("M132").Value2 = -127543.121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value2 = 669989
$ws.Range("J20").Value2 = 9969
$ws.Range("L20").Value2 = 9969
$ws.Range("N20").Value2 = -10421
$ws.Range("H23").Value2 = 19750
$ws.Range("I23").Value2 = 19750
$ws.Range("K23").Value2 = 19750
$ws.Range("M23").Value2 = -19520
$ws.Range("H29").Value2 = 59999
$ws.Range("I29").Value2 = 39999
$ws.Range("J29").Value2 = 99999
$ws.Range("K29").Value2 = 39999
$ws.Range("L29").Value2 = 99999
$ws.Range("M29").Value2 = -39704
$ws.Range("N29").Value2 = -100589
$ws.Range("H35").Value2 = 2934.4443
$ws.Range("I35").Value2 = 2934.4443
$ws.Range("J35").Value2 = 0
$ws.Range("K35").Value2 = 2934.4443
$ws.Range("L35").Value2 = 0
$ws.Range("N35").Value2 = -2598.4443
$ws.Range("H61").Value2 = 1974.3334
$ws.Range("I61").Value2 = 1949.2727
$ws.Range("K61").Value2 = 1949.2727
$ws.Range("M61").Value2 = -1747.2727
$ws.Range("H62").Value2 = 36662
$ws.Range("J62").Value2 = 36662
$ws.Range("L62").Value2 = 36662
$ws.Range("N62").Value2 = -37910
$ws.Range("H65").Value2 = 36662
$ws.Range("J65").Value2 = 36662
$ws.Range("L65").Value2 = 109986
$ws.Range("N65").Value2 = -116226
$ws.Range("H113").Value2 = 1974.3334
$ws.Range("I113").Value2 = 1949.2727
$ws.Range("K113").Value2 = 1949.2727
$ws.Range("M113").Value2 = 220.7273
$ws.Range("H122").Value2 = 7411.0625
$ws.Range("I122").Value2 = 8712.714
$ws.Range("K122").Value2 = 26138.142
$ws.Range("M122").Value2 = -23688.142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value2 = 18700
$ws.Range("I32").Value2 = 8400
$ws.Range("K32").Value2 = 8400
$ws.Range("M32").Value2 = -8083
$ws.Range("H81").Value2 = 2169.8235
$ws.Range("I81").Value2 = 2246.4
$ws.Range("J81").Value2 = 1595.5
$ws.Range("K81").Value2 = 4492.8
$ws.Range("L81").Value2 = 3191
$ws.Range("M81").Value2 = -3431.8
$ws.Range("N81").Value2 = -5313
$ws.Range("H84").Value2 = 2169.8235
$ws.Range("I84").Value2 = 2246.4
$ws.Range("J84").Value2 = 1595.5
$ws.Range("K84").Value2 = 22464
$ws.Range("L84").Value2 = 15955
$ws.Range("M84").Value2 = -17160
$ws.Range("N84").Value2 = -26563
$ws.Range("H100").Value2 = 737.0454999999999
$ws.Range("I100").Value2 = 641.15
$ws.Range("J100").Value2 = 1696
$ws.Range("K100").Value2 = 1282.3
$ws.Range("L100").Value2 = 3392
$ws.Range("M100").Value2 = -741.3
$ws.Range("N100").Value2 = -4474
$ws.Range("H122").Value2 = 2368.4614
$ws.Range("I122").Value2 = 2274.1667
$ws.Range("K122").Value2 = 6822.500100000001
$ws.Range("M122").Value2 = -4372.500100000001
$ws.Range("H126").Value2 = 2305.1875
$ws.Range("I126").Value2 = 2157.4167
$ws.Range("K126").Value2 = 6472.250100000001
$ws.Range("M126").Value2 = -4002.250100000001
